$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Junio de 2020 a las 19:19"

# Swap country name labels that changed rank/order
$tmpA = $ws.Range("A98").Value2
$tmpB = $ws.Range("A99").Value2
$ws.Range("A98").Value = $tmpB
$ws.Range("A99").Value = $tmpA

$tmpA = $ws.Range("A145").Value2
$tmpB = $ws.Range("A146").Value2
$ws.Range("A145").Value = $tmpB
$ws.Range("A146").Value = $tmpA

$tmpA = $ws.Range("A202").Value2
$tmpB = $ws.Range("A203").Value2
$ws.Range("A202").Value = $tmpB
$ws.Range("A203").Value = $tmpA

$tmpA = $ws.Range("A208").Value2
$tmpB = $ws.Range("A209").Value2
$ws.Range("A208").Value = $tmpB
$ws.Range("A209").Value = $tmpA

$tmpA = $ws.Range("A211").Value2
$tmpB = $ws.Range("A212").Value2
$ws.Range("A211").Value = $tmpB
$ws.Range("A212").Value = $tmpA

# Update numeric statistics cells
$ws.Range("B4").Value = 2479551
$ws.Range("C4").Value = 16997
$ws.Range("D4").Value = 1041307
$ws.Range("E4").Value = 1313726
$ws.Range("G4").Value = 237
$ws.Range("H4").Value = 124518

$ws.Range("B5").Value = 1207721
$ws.Range("C5").Value = 15247
$ws.Range("E5").Value = 503379
$ws.Range("G5").Value = 560
$ws.Range("H5").Value = 54434

$ws.Range("B7").Value = 489960
$ws.Range("C7").Value = 16975
$ws.Range("D7").Value = 285211
$ws.Range("E7").Value = 189446
$ws.Range("G7").Value = 396
$ws.Range("H7").Value = 15303

$ws.Range("B11").Value = 259064
$ws.Range("C11").Value = 4648
$ws.Range("D11").Value = 219327
$ws.Range("E11").Value = 34834
$ws.Range("G11").Value = 172
$ws.Range("H11").Value = 4903

$ws.Range("B12").Value = 239706
$ws.Range("C12").Value = 296
$ws.Range("D12").Value = 186725
$ws.Range("E12").Value = 18303
$ws.Range("G12").Value = 34
$ws.Range("H12").Value = 34678

$ws.Range("B15").Value = 193299
$ws.Range("C15").Value = 45
$ws.Range("E15").Value = 7495
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 9004

$ws.Range("B48").Value = 25405
$ws.Range("C48").Value = 9
$ws.Range("E48").Value = 314
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 1727

$ws.Range("B52").Value = 22400
$ws.Range("C52").Value = 356
$ws.Range("D52").Value = 16007
$ws.Range("E52").Value = 6084
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 309

$ws.Range("B66").Value = 12445
$ws.Range("C66").Value = 197
$ws.Range("D66").Value = 8920
$ws.Range("E66").Value = 2647
$ws.Range("G66").Value = 9
$ws.Range("H66").Value = 878

$ws.Range("B98").Value = 3033
$ws.Range("C98").Value = 206
$ws.Range("D98").Value = 1096
$ws.Range("E98").Value = 1927
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 10

$ws.Range("B99").Value = 2878
$ws.Range("C99").Value = 43
$ws.Range("D99").Value = 868
$ws.Range("E99").Value = 1920
$ws.Range("H99").Value = 90

$ws.Range("B133").Value = 992
$ws.Range("C133").Value = 1
$ws.Range("E133").Value = 149

$ws.Range("B138").Value = 863
$ws.Range("C138").Value = 3
$ws.Range("D138").Value = 774
$ws.Range("E138").Value = 15

$ws.Range("B141").Value = 821
$ws.Range("C141").Value = 16
$ws.Range("D141").Value = 731
$ws.Range("E141").Value = 90

$ws.Range("B145").Value = 706
$ws.Range("C145").Value = 16
$ws.Range("D145").Value = 347
$ws.Range("E145").Value = 351
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 8

$ws.Range("B146").Value = 698
$ws.Range("D146").Value = 653
$ws.Range("E146").Value = 3
$ws.Range("H146").Value = 42

$ws.Range("B156").Value = 411
$ws.Range("C156").Value = 22
$ws.Range("E156").Value = 87

$ws.Range("D161").Value = 211
$ws.Range("E161").Value = 76

$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0
